# Update the "Price" (D) and "Volume(1h)" (E) columns of the cryptos sheet
# with freshly scraped values. Price values are forced to text (leading
# apostrophe) so that numeric-looking strings (e.g. "0.999", "2.80",
# "37.210.52") are kept as literal text instead of being coerced to
# Excel numbers, matching the original inline-string cell contents.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''37.210.52'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '''2.061.70'
$ws.Range('E3').Value = '  -1.49%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '''250.14'
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range('D6').Value = '''0.677'
$ws.Range('E6').Value = '  +2.91%  '
$ws.Range('D7').Value = '''59.66'
$ws.Range('E7').Value = '  +17.90%  '
$ws.Range('D8').Value = '''0.999'
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').Value = '''60.73'
$ws.Range('E9').Value = '  -0.51%  '
$ws.Range('D10').Value = '''0.379'
$ws.Range('E10').Value = '  +1.48%  '
$ws.Range('D11').Value = '''0.0802'
$ws.Range('E11').Value = '  +7.31%  '
$ws.Range('E12').Value = '  +1.20%  '
$ws.Range('D13').Value = '''15.27'
$ws.Range('E13').Value = '  +0.78%  '
$ws.Range('D14').Value = '''2.358.42'
$ws.Range('E14').Value = '  -1.28%  '
$ws.Range('D15').Value = '''0.819'
$ws.Range('E15').Value = '  -2.16%  '
$ws.Range('D16').Value = '''5.34'
$ws.Range('E16').Value = '  +4.26%  '
$ws.Range('D17').Value = '''2.060.39'
$ws.Range('E17').Value = '  -1.37%  '
$ws.Range('D18').Value = '''37.132.47'
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('D19').Value = '''75.21'
$ws.Range('E19').Value = '  +3.59%  '
$ws.Range('D20').Value = '''0.0₃0921'
$ws.Range('E20').Value = '  +11.87%  '
$ws.Range('D21').Value = '''14.58'
$ws.Range('E21').Value = '  +8.98%  '
$ws.Range('D22').Value = '''5.37'
$ws.Range('E22').Value = '  +2.46%  '
$ws.Range('D23').Value = '''239.26'
$ws.Range('E23').Value = '  -0.62%  '
$ws.Range('E24').Value = '  +0.07%  '
$ws.Range('D25').Value = '''2.45'
$ws.Range('E25').Value = '  -2.00%  '
$ws.Range('D26').Value = '''171.93'
$ws.Range('E26').Value = '  +1.19%  '
$ws.Range('D27').Value = '''9.19'
$ws.Range('E27').Value = '  -1.50%  '
$ws.Range('D28').Value = '''20.42'
$ws.Range('E28').Value = '  -4.54%  '
$ws.Range('E29').Value = '  -0.34%  '
$ws.Range('E30').Value = '  +2.37%  '
$ws.Range('D31').Value = '''4.63'
$ws.Range('E31').Value = '  +2.57%  '
$ws.Range('D32').Value = '''1.08'
$ws.Range('E32').Value = '  -3.02%  '
$ws.Range('D33').Value = '''0.0634'
$ws.Range('E33').Value = '  +4.08%  '
$ws.Range('D34').Value = '''4.42'
$ws.Range('E34').Value = '  +7.64%  '
$ws.Range('D35').Value = '''0.0886'
$ws.Range('E35').Value = '  -4.76%  '
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('D37').Value = '''2.27'
$ws.Range('E37').Value = '  -1.07%  '
$ws.Range('E38').Value = '  -2.81%  '
$ws.Range('D39').Value = '''0.110'
$ws.Range('E39').Value = '  +27.31%  '
$ws.Range('E40').Value = '  +1.57%  '
$ws.Range('D41').Value = '''18.41'
$ws.Range('E41').Value = '  +4.68%  '
$ws.Range('E42').Value = '  +0.64%  '
$ws.Range('E43').Value = '  -0.30%  '
$ws.Range('D44').Value = '''4.41'
$ws.Range('E44').Value = '  +31.06%  '
$ws.Range('D45').Value = '''97.48'
$ws.Range('E45').Value = '  -0.53%  '
$ws.Range('E46').Value = '  -0.21%  '
$ws.Range('D47').Value = '''4.47'
$ws.Range('E47').Value = '  +13.15%  '
$ws.Range('D48').Value = '''2.51'
$ws.Range('E48').Value = '  +10.47%  '
$ws.Range('D49').Value = '''1.304.62'
$ws.Range('E49').Value = '  -1.26%  '
$ws.Range('D50').Value = '''2.92'
$ws.Range('E50').Value = '  -1.68%  '
$ws.Range('D51').Value = '''6.90'
$ws.Range('E51').Value = '  -0.29%  '
